$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scale the values in column D (rows 2-33) by 10,000 (push/pull data scale change).
# Literal target values are used (rather than a live *10000 multiplication) so the
# stored doubles match exactly bit-for-bit, avoiding floating point drift.
$values = @{
    2  = 1696693.946436
    3  = 11475286.204935
    4  = 3675993.415221
    5  = 1181036.168182
    6  = 6521066.011809
    7  = 576560.29715
    8  = 1761018.289914
    9  = 281751.117122
    10 = 1453071.446718
    11 = 1276013.107121
    12 = 1076757.220732
    13 = 4949463.434506
    14 = 1516117.153227
    15 = 545214.22245
    16 = 10163395.752203
    17 = 1077100.14431
    18 = 2312704.725347
    19 = 1885988.776212
    20 = 286118.210905
    21 = 471048.657485
    22 = 98357896
    23 = 3865607.331231
    24 = 972416.554906
    25 = 2985237.182962
    26 = 5101282.045263
    27 = 2932643.595974
    28 = 3307078.645489
    29 = 169821.278996
    30 = 1212578.103854
    31 = 1735785.95345
    32 = 2019061.128324
    33 = 1843372.885213
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 4).Value2 = $values[$r]
}
